# Adds a new "2022-Q1" sheet (fund-holdings detail, like the other
# quarterly sheets) positioned right before "总计", and updates the
# "总计" (totals) roll-up sheet with a new leading row for 2022-Q1.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Insert the new "2022-Q1" worksheet before "总计"
# ---------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")
$qSheet = $wb.Worksheets.Add($totalSheet)
$qSheet.Name = "2022-Q1"

# Re-resolve "总计" by name: the handle captured above goes stale once
# a new sheet is spliced in next to it (it ends up aliasing the newly
# added sheet instead), so grab a fresh, live reference to it here.
$totalSheet = $wb.Worksheets.Item("总计")

# Use "2021-Q4" as the style/template donor for the header row and the
# "2021-Q3" sheet (24 data rows) as the donor for the bulk of the data
# rows so the new sheet inherits the same look (bold+bordered header,
# bold index column) as its siblings.
$hdrSrc = $wb.Worksheets.Item("2021-Q4")
$rowSrc = $wb.Worksheets.Item("2021-Q3")

$hdrSrc.Range("B1:H1").Copy($qSheet.Range("B1:H1"))
$rowSrc.Range("A2:H25").Copy($qSheet.Range("A2:H25"))
for ($r = 26; $r -le 28; $r++) {
    $rowSrc.Range("A25:H25").Copy($qSheet.Range("A" + $r + ":H" + $r))
}

$qSheet.Range("B1").Value = "基金代码"
$qSheet.Range("C1").Value = "基金名称"
$qSheet.Range("D1").Value = "基金规模"
$qSheet.Range("E1").Value = "股票总仓位"
$qSheet.Range("F1").Value = "仓位占比"
$qSheet.Range("G1").Value = "持有市值(亿元)"
$qSheet.Range("H1").Value = "仓位排名"

$fundRows = @(
    @("519087", "新华优选分红混合", "9.03", "88.56", "7.92", "0.7152", 1),
    @("001040", "新华策略精选股票", "6.15", "93.72", "8.09", "0.4975", 1),
    @("519156", "新华行业轮换灵活配置混合A", "4.98", "93.77", "7.85", "0.3909", 1),
    @("003291", "信达澳银健康中国灵活配置混合", "4.03", "91.32", "5.00", "0.2015", 7),
    @("011598", "信达澳银医药健康混合", "3.47", "91.98", "4.75", "0.1648", 7),
    @("014185", "招商专精特新股票A", "8.37", "30.94", "1.87", "0.1565", 4),
    @("001294", "新华战略新兴产业灵活配置混合", "1.07", "93.41", "8.29", "0.0887", 1),
    @("005108", "圆信永丰双利优选定期开放灵活配置混合", "1.89", "94.60", "4.27", "0.0807", 7),
    @("011383", "富安达医药创新混合", "1.68", "83.50", "4.61", "0.0774", 3),
    @("013067", "富安达中小盘六个月持有期混合", "2.45", "74.39", "2.97", "0.0728", 7),
    @("006981", "中金新医药股票A", "1.86", "91.77", "3.87", "0.0720", 5),
    @("011457", "新华行业龙头主题股票", "0.85", "93.55", "7.91", "0.0672", 1),
    @("014186", "招商专精特新股票C", "3.46", "30.94", "1.87", "0.0647", 4),
    @("519097", "新华中小市值优选混合", "0.75", "62.70", "5.01", "0.0376", 3),
    @("001965", "圆信永丰兴源灵活配置混合A", "0.76", "93.43", "4.60", "0.0350", 8),
    @("001861", "富安达健康人生灵活配置混合", "0.61", "82.18", "5.13", "0.0313", 2),
    @("007861", "金元顺安医疗健康混合型证券投资基金A", "0.52", "86.80", "3.31", "0.0172", 10),
    @("007005", "中金新医药股票C", "0.34", "91.77", "3.87", "0.0132", 5),
    @("001966", "圆信永丰兴源灵活配置混合C", "0.25", "93.43", "4.60", "0.0115", 8),
    @("008884", "博远博锐混合A", "0.19", "86.59", "4.62", "0.0088", 7),
    @("006274", "圆信永丰医药健康混合", "0.18", "93.66", "4.63", "0.0083", 8),
    @("519157", "新华行业轮换灵活配置混合C", "0.04", "93.77", "7.85", "0.0031", 1),
    @("007862", "金元顺安医疗健康混合型证券投资基金C", "0.09", "86.80", "3.31", "0.0030", 10),
    @("001659", "富安达新动力灵活配置混合", "0.07", "90.52", "3.91", "0.0027", 10),
    @("005537", "中航新起航灵活配置混合A", "0.03", "87.09", "4.85", "0.0015", 8),
    @("008885", "博远博锐混合C", "0.02", "86.59", "4.62", "0.0009", 7),
    @("005538", "中航新起航灵活配置混合C", "0.01", "87.09", "4.85", "0.0005", 8)
)

$r = 2
foreach ($row in $fundRows) {
    $qSheet.Range("A" + $r).Value = ($r - 2)
    $qSheet.Range("B" + $r).Value = "'" + $row[0]
    $qSheet.Range("C" + $r).Value = $row[1]
    $qSheet.Range("D" + $r).Value = "'" + $row[2]
    $qSheet.Range("E" + $r).Value = "'" + $row[3]
    $qSheet.Range("F" + $r).Value = "'" + $row[4]
    $qSheet.Range("G" + $r).Value = "'" + $row[5]
    $qSheet.Range("H" + $r).Value = $row[6]
    $r++
}

# ---------------------------------------------------------------------
# 2. Prepend a "2022-Q1" row to the "总计" roll-up sheet
# ---------------------------------------------------------------------
$totalSheet.Rows.Item(2).Insert()

# The blank inserted row borrows row 1's (header) bold/centered style;
# restore column A's normal index-column look and clear B:D back to the
# default (unstyled) look used by every other data row.
$totalSheet.Range("A3").Copy($totalSheet.Range("A2"))
$totalSheet.Range("B3:D3").Copy($totalSheet.Range("B2:D2"))

$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q1"
$totalSheet.Range("C2").Value = 27
$totalSheet.Range("D2").Value = 2.82

# Renumber the index column for the rows that shifted down one spot.
$totalSheet.Range("A3").Value = 1
$totalSheet.Range("A4").Value = 2
$totalSheet.Range("A5").Value = 3
$totalSheet.Range("A6").Value = 4
